$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 6 (Zeit = 120), shifting rows 7-9 up to 6-8
$ws.Rows.Item(6).Delete()

# Update selection to match the final state observed after the edit
$ws.Range("F13").Select()
